# Michael_Schulze.xlsx - "USB2 read daten anstatt USB3"
# Replace column A readings with the USB2 measurement series, and format the
# (until now unused) column B with a 3-decimal number format, matching the
# newly added B1:B10 cells that accompany the reading column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(17.387, 16.96, 17.639, 16.648, 18.541, 18.024, 18.553, 17.675, 18.169, 16.747)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
    $ws.Cells.Item($row, 2).NumberFormat = "0.000"
}

# Header/footer font label changed from "Standard" to "Regular"
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Seite &P'

# Selection moved to A11 (first empty row after the data) after the edit
$ws.Range("A11").Select()
